$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "municipio-nombre" column (E) metadata rows are re-curated from a
# measure into a dimension, matching the shape already used by the
# "provincia-nombre" column (F):
#   E2 (sdmx type) : iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   E3 (dim/medida): medida                          -> dim
#   E4 (datatype)  : xsd:int                          -> URI-Municipio
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"
